$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, D (fecha/serial), L (calidad), M (volumen), N/O/P (precios),
# Q (unidad comercializacion), R (origen), S (precio $/Kg), T (Kg/unidad)
$rows = @(
    [pscustomobject]@{ Row=2; D=44868; L="Primera"; M=30; N=14000; O=14000; P=14000; Q="$/bandeja 5 kilos"; R="Provincia de Quillota"; S=2800; T=5 }
    [pscustomobject]@{ Row=3; D=44519; L="Primera"; M=30; N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=2800; T=10 }
    [pscustomobject]@{ Row=4; D=44496; L="Primera"; M=55; N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=2800; T=10 }
    [pscustomobject]@{ Row=5; D=44889; L="Primera"; M=50; N=30000; O=30000; P=30000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=3000; T=10 }
    [pscustomobject]@{ Row=6; D=44879; L="Primera"; M=25; N=30000; O=30000; P=30000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=3000; T=10 }
    [pscustomobject]@{ Row=7; D=44483; L="Primera"; M=35; N=10000; O=10000; P=10000; Q="$/bandeja 5 kilos"; R="Provincia de Quillota"; S=2000; T=5 }
    [pscustomobject]@{ Row=8; D=44874; L="Primera"; M=40; N=25000; O=25000; P=25000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=2500; T=10 }
    [pscustomobject]@{ Row=9; D=44902; L="Primera"; M=90; N=25000; O=25000; P=25000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=2500; T=10 }
    [pscustomobject]@{ Row=10; D=44466; L="Primera"; M=80; N=11000; O=11000; P=11000; Q="$/bandeja 5 kilos"; R="La Ligua"; S=2200; T=5 }
    [pscustomobject]@{ Row=11; D=44515; L="Primera"; M=80; N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Los Andes"; S=2800; T=10 }
    [pscustomobject]@{ Row=12; D=44488; L="Primera"; M=100; N=12000; O=12000; P=12000; Q="$/bandeja 5 kilos"; R="La Ligua"; S=2400; T=5 }
    [pscustomobject]@{ Row=13; D=44859; L="Primera"; M=30; N=20000; O=20000; P=20000; Q="$/bandeja 5 kilos"; R="Provincia de Quillota"; S=4000; T=5 }
    [pscustomobject]@{ Row=14; D=44921; L="Primera"; M=55; N=15000; O=15000; P=15000; Q="$/bandeja 7 kilos"; R="Provincia de Quillota"; S=2143; T=7 }
    [pscustomobject]@{ Row=15; D=44858; L="Primera"; M=90; N=20000; O=20000; P=20000; Q="$/bandeja 5 kilos"; R="Provincia de Quillota"; S=4000; T=5 }
    [pscustomobject]@{ Row=16; D=44901; L="Primera"; M=40; N=25000; O=25000; P=25000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=2500; T=10 }
    [pscustomobject]@{ Row=17; D=44511; L="Primera"; M=45; N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Los Andes"; S=2800; T=10 }
    [pscustomobject]@{ Row=18; D=44511; L="Primera"; M=45; N=3200; O=3200; P=3200; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=320; T=10 }
    [pscustomobject]@{ Row=19; D=44503; L="Primera"; M=50; N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=2800; T=10 }
    [pscustomobject]@{ Row=20; D=44166; L="Segunda"; M=20; N=12000; O=12000; P=12000; Q="$/caja 18 kilos"; R="La Ligua"; S=667; T=18 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N: Precio minimo
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Precio maximo
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 17).Value = $r.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r.Row, 18).Value = $r.R   # R: Origen
    $ws.Cells.Item($r.Row, 19).Value = $r.S   # S: Precio $/Kg
    $ws.Cells.Item($r.Row, 20).Value = $r.T   # T: Kg / unidad
}
